$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a value to be stored as Text, matching the original
# inlineStr cell type (prevents Excel's automatic number/date inference
# from collapsing things like '233.51' into a float or dropping the
# trailing zero in values like '58.30'), while leaving the cell style
# index unchanged (no new applied number format remains on the cell).
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '37.748.76'
$ws.Range("E2").Value = '  +0.08%  '
Set-TextValue "D3" '2.077.28'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue "D5" '233.51'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  +0.02%  '
Set-TextValue "D8" '58.30'
$ws.Range("E8").Value = '  +0.82%  '
Set-TextValue "D9" '0.394'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("E11").Value = '  +2.65%  '
Set-TextValue "D12" '2.382.42'
$ws.Range("E12").Value = '  -1.27%  '
Set-TextValue "D13" '14.76'
$ws.Range("E13").Value = '  +1.85%  '
Set-TextValue "D14" '20.87'
$ws.Range("E14").Value = '  -2.24%  '
Set-TextValue "D15" '0.775'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("E16").Value = '  +2.27%  '
Set-TextValue "D17" '2.111.63'
$ws.Range("E17").Value = '  +0.20%  '
Set-TextValue "D18" '37.695.41'
$ws.Range("E18").Value = '  +0.00%  '
Set-TextValue "D19" '6.16'
$ws.Range("E19").Value = '  -0.62%  '
Set-TextValue "D20" '71.02'
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("E21").Value = '  +1.29%  '
Set-TextValue "D22" '228.24'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("E23").Value = '  +0.00%  '
Set-TextValue "D24" '2.38'
$ws.Range("E24").Value = '  -0.81%  '
Set-TextValue "D25" '2.40'
$ws.Range("E25").Value = '  -0.35%  '
Set-TextValue "D26" '170.67'
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("E27").Value = '  +5.49%  '
Set-TextValue "D28" '9.03'
$ws.Range("E28").Value = '  +1.36%  '
Set-TextValue "D29" '19.45'
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("E31").Value = '  +2.72%  '
Set-TextValue "D32" '4.69'
$ws.Range("E32").Value = '  +1.50%  '
Set-TextValue "D33" '0.0629'
$ws.Range("E33").Value = '  +1.29%  '
Set-TextValue "D34" '4.65'
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("E35").Value = '  -3.88%  '
Set-TextValue "D36" '1.82'
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("E37").Value = '  -1.57%  '
Set-TextValue "D38" '0.999'
$ws.Range("E38").Value = '  -0.30%  '
Set-TextValue "D39" '5.30'
$ws.Range("E39").Value = '  -1.23%  '
Set-TextValue "D40" '0.0971'
$ws.Range("E40").Value = '  -1.82%  '
Set-TextValue "D41" '98.56'
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("E42").Value = '  -2.29%  '
Set-TextValue "D43" '0.0214'
$ws.Range("E43").Value = '  +0.89%  '
Set-TextValue "D44" '1.449.49'
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D45" '1.15'
$ws.Range("E45").Value = '  -1.50%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D46" '16.38'
$ws.Range("E46").Value = '  +6.67%  '
Set-TextValue "D47" '4.20'
$ws.Range("E47").Value = '  +1.79%  '
$ws.Range("E48").Value = '  +1.22%  '
Set-TextValue "D49" '7.39'
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("E50").Value = '  -0.65%  '
Set-TextValue "D51" '2.267.54'
$ws.Range("E51").Value = '  -1.38%  '
